$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spotten")

# New scenario row: C2 gets "Yes", formatted like the other green "done" cells (e.g. C12)
$ws.Range("C2").Value = "Yes"
$ws.Range("C12").Copy()
$ws.Range("C2").PasteSpecial(-4122)  # xlPasteFormats

# Scroll/selection state: no longer pinned at row 7, now selecting C4
$ws.Range("C4").Select()

# Window was resized (un-maximized) by the author
$excel.ActiveWindow.Width = 15247
$excel.ActiveWindow.Height = 11647
